$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the "Units" column (E) across the first (ADNI_EMBICDCB / UPENNPLASMA) block ---
$ws.Range("E1:E13").ClearContents()

# --- Update the UPENNPLASMA rows (12-13): drop per-row units, keep description text ---
$ws.Range("C12").Value = "Abeta1-40 result in plasma"
$ws.Range("C13").Value = "Abeta1-42 result in plasma"

# --- Wipe everything from row 14 through row 23 (old "Not Time Bound" block lived here,
#     now relocating down to make room for the new UPENNBIOMK_MASTER rows) ---
$ws.Range("A14:I23").ClearContents()

# --- New UPENNBIOMK_MASTER rows (14-19) ---
$ws.Range("A14").Value = "UPENNBIOMK_MASTER"
$ws.Range("B14").Value = "ABETA"
$ws.Range("C14").Value = "Normalized Amyloid Beta (1-42) in CSF"
$ws.Range("D14").Value = "Ratio"

$ws.Range("A15").Value = "UPENNBIOMK_MASTER"
$ws.Range("B15").Value = "TAU"
$ws.Range("C15").Value = "Normalized TAU in CSF"
$ws.Range("D15").Value = "Ratio"

$ws.Range("A16").Value = "UPENNBIOMK_MASTER"
$ws.Range("B16").Value = "PTAU"
$ws.Range("C16").Value = "Normalized PTAU (181) in CSF"
$ws.Range("D16").Value = "Ratio"

$ws.Range("A17").Value = "UPENNBIOMK_MASTER"
$ws.Range("B17").Value = "ABETA_RAW"
$ws.Range("C17").Value = "Raw amyloid beta (1-42) in CSF"
$ws.Range("D17").Value = "Ratio"

$ws.Range("A18").Value = "UPENNBIOMK_MASTER"
$ws.Range("B18").Value = "TAU_RAW"
$ws.Range("C18").Value = "Raw tau in CSF"
$ws.Range("D18").Value = "Ratio"

$ws.Range("A19").Value = "UPENNBIOMK_MASTER"
$ws.Range("B19").Value = "PTAU_RAW"
$ws.Range("C19").Value = "Raw ptau (181) beta in CSF"
$ws.Range("D19").Value = "Ratio"

# --- "Not Time Bound" section header + methods note, now at row 23 (was row 16) ---
$ws.Range("A23").Value = "Not Time Bound"
$ws.Range("G23").Value = "See ""ADNI_Biomarker_Methods_Statistical_Analyses_Dec2023.docx"" for info on Z-score calculations"

# --- Column header row, now at row 24 (was row 17); Units column dropped ---
$ws.Range("A24").Value = "Source Table"
$ws.Range("B24").Value = "Variable"
$ws.Range("C24").Value = "Description"
$ws.Range("D24").Value = "Type"

# --- Wipe the old ADSP_PHC_CSF_Dec2023 block (rows 18-23 previously) which now lives at 25-30 ---
$ws.Range("A25:E30").ClearContents()

# --- ADSP_PHC_CSF block, now rows 25-30 (was rows 18-23); Units column dropped, table renamed ---
$ws.Range("A25").Value = "ADSP_PHC_CSF"
$ws.Range("B25").Value = "AB42_RAW"
$ws.Range("C25").Value = "Raw AB42 biomarker levels"
$ws.Range("D25").Value = "Ratio"

$ws.Range("A26").Value = "ADSP_PHC_CSF"
$ws.Range("B26").Value = "PHC_AB42"
$ws.Range("C26").Value = "Harmonized biomarker AB42 Z-score (derived - see METHODS)"
$ws.Range("D26").Value = "Ratio"

$ws.Range("A27").Value = "ADSP_PHC_CSF"
$ws.Range("B27").Value = "Tau_RAW"
$ws.Range("C27").Value = "Raw Tau biomarker levels"
$ws.Range("D27").Value = "Ratio"

$ws.Range("A28").Value = "ADSP_PHC_CSF"
$ws.Range("B28").Value = "PHC_Tau"
$ws.Range("C28").Value = "Harmonized biomarker Tau Z-score (derived - see METHODS)"
$ws.Range("D28").Value = "Ratio"

$ws.Range("A29").Value = "ADSP_PHC_CSF"
$ws.Range("B29").Value = "pTau_RAW"
$ws.Range("C29").Value = "Raw pTau biomarker levels"
$ws.Range("D29").Value = "Ratio"

$ws.Range("A30").Value = "ADSP_PHC_CSF"
$ws.Range("B30").Value = "PHC_pTau"
$ws.Range("C30").Value = "Harmonized biomarker pTau Z-score (derived - see METHODS)"
$ws.Range("D30").Value = "Ratio"

# --- New SELKOELAB_OAB rows (31-33) ---
$ws.Range("A31").Value = "SELKOELAB_OAB"
$ws.Range("B31").Value = "MEAN"
$ws.Range("C31").Value = "Mean amount of A-Beta Oligomer"
$ws.Range("D31").Value = "Ratio"

$ws.Range("A32").Value = "SELKOELAB_OAB"
$ws.Range("B32").Value = "SD"
$ws.Range("C32").Value = "Standard deviation"
$ws.Range("D32").Value = "Ratio"

$ws.Range("A33").Value = "SELKOELAB_OAB"
$ws.Range("B33").Value = "CV"
$ws.Range("C33").Value = "Coefficient of Variation"
$ws.Range("D33").Value = "Ratio"

# --- Column A width (best effort) ---
$ws.Columns.Item(1).ColumnWidth = 14

# --- Sheet view: reposition selection to I12 (also clears the old topLeftCell scroll anchor) ---
$ws.Range("I12").Select()
